$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.445.49'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '2.019.17'
$ws.Range('E3').Value = '  +6.42%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''246.34'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').Value = '''0.661'
$ws.Range('E6').Value = '  -4.64%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '''45.18'
$ws.Range('E8').Value = '  +4.94%  '
$ws.Range('D9').Value = '''59.64'
$ws.Range('E9').Value = '  +5.64%  '
$ws.Range('D10').Value = '''0.367'
$ws.Range('E10').Value = '  +2.64%  '
$ws.Range('D11').Value = '''0.0717'
$ws.Range('E11').Value = '  -5.45%  '
$ws.Range('D12').Value = '''0.0985'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '''14.60'
$ws.Range('E13').Value = '  +3.64%  '
$ws.Range('D14').Value = '2.315.83'
$ws.Range('E14').Value = '  +6.81%  '
$ws.Range('D15').Value = '''0.809'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '2.023.06'
$ws.Range('E16').Value = '  +6.31%  '
$ws.Range('D17').Value = '''4.92'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('D18').Value = '36.394.59'
$ws.Range('E18').Value = '  +2.79%  '
$ws.Range('D19').Value = '''71.38'
$ws.Range('E19').Value = '  -3.15%  '
$ws.Range('D20').Value = '0.0₃0821'
$ws.Range('E20').Value = '  -1.51%  '
$ws.Range('D21').Value = '''12.99'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '''236.47'
$ws.Range('E22').Value = '  -3.41%  '
$ws.Range('D23').Value = '''4.89'
$ws.Range('E23').Value = '  -6.44%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  -9.03%  '
$ws.Range('D26').Value = '''163.40'
$ws.Range('E26').Value = '  -2.46%  '
$ws.Range('D27').Value = '''19.72'
$ws.Range('E27').Value = '  +7.42%  '
$ws.Range('D28').Value = '''8.61'
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('D29').Value = '''1.95'
$ws.Range('E29').Value = '  -10.42%  '
$ws.Range('D30').Value = '''0.122'
$ws.Range('E30').Value = '  -4.64%  '
$ws.Range('D31').Value = '''22.36'
$ws.Range('E31').Value = '  +61.96%  '
$ws.Range('D32').Value = '''4.42'
$ws.Range('E32').Value = '  +1.49%  '
$ws.Range('D33').Value = '''0.0591'
$ws.Range('E33').Value = '  -2.31%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D36').Value = '''4.01'
$ws.Range('E36').Value = '  -5.96%  '
$ws.Range('D37').Value = '''0.0810'
$ws.Range('E37').Value = '  +10.01%  '
$ws.Range('D38').Value = '''2.13'
$ws.Range('E38').Value = '  +8.24%  '
$ws.Range('D39').Value = '''0.847'
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('D40').Value = '''1.34'
$ws.Range('E40').Value = '  -9.36%  '
$ws.Range('D41').Value = '''0.0216'
$ws.Range('E41').Value = '  -4.20%  '
$ws.Range('D42').Value = '''96.07'
$ws.Range('E42').Value = '  -3.24%  '
$ws.Range('D43').Value = '''1.11'
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').Value = '''2.75'
$ws.Range('E44').Value = '  +13.97%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '''16.02'
$ws.Range('E45').Value = '  -5.65%  '
$ws.Range('D46').Value = '1.320.44'
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('D47').Value = '''0.0814'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').Value = '''2.77'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''2.22'
$ws.Range('E49').Value = '  -6.53%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.196.78'
$ws.Range('E50').Value = '  +6.17%  '
$ws.Range('D51').Value = '''3.84'
$ws.Range('E51').Value = '  +15.25%  '
